# Common: Added some initial mix detail
# Appends new translation rows (lab.mixture.index.* / lab.mixture.preview.*)
# to the "Import" sheet, right after the existing last row (495).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# New rows keep the same formatting (style) as the preceding data rows,
# so copy the format down first.
$ws.Range("A495:C495").Copy()
$ws.Range("A496:C503").PasteSpecial(-4122)  # xlPasteFormats

$data = @(
    @("cs", "lab.mixture.index.title", "Detail mixu"),
    @("cs", "lab.mixture.index.preview.title", "Náhled mixu"),
    @("cs", "lab.mixture.index.preview.subtitle", "Zde můžete spravovat vybraný mix."),
    @("cs", "lab.mixture.preview.name", "Název"),
    @("cs", "lab.mixture.preview.liquid", "Liquid"),
    @("cs", "lab.mixture.preview.nicotine", "Obsah nikotinu"),
    @("cs", "lab.mixture.preview.base", "Báze"),
    @("cs", "lab.mixture.preview.booster", "Booster")
)

$row = 496
foreach ($entry in $data) {
    $ws.Range("A$row").Value = $entry[0]
    $ws.Range("B$row").Value = $entry[1]
    $ws.Range("C$row").Value = $entry[2]
    $row++
}

# Match the author's final cursor / scroll position.
$ws.Range("B498").Select()
